# Scheduled-runner update: refresh cached market-board price figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across the
# per-job leve-profit sheets, row by row, matching the upstream export.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 71435820
$ws.Range("I74").Value = 136366500
$ws.Range("K74").Value = 136366500
$ws.Range("M74").Value = -136365564
$ws.Range("H77").Value = 71435820
$ws.Range("I77").Value = 136366500
$ws.Range("K77").Value = 681832500
$ws.Range("M77").Value = -681827820
$ws.Range("H112").Value = 4784.1665
$ws.Range("J112").Value = 4784.1665
$ws.Range("L112").Value = 14352.4995
$ws.Range("N112").Value = -16568.4995
$ws.Range("H132").Value = 1856.3778
$ws.Range("I132").Value = 1758.641
$ws.Range("K132").Value = 5275.923000000001
$ws.Range("M132").Value = -2745.923000000001
$ws.Range("H137").Value = 3594.4443
$ws.Range("I137").Value = 3916.8125
$ws.Range("K137").Value = 11750.4375
$ws.Range("M137").Value = -9200.4375

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1528695
$ws.Range("I32").Value = 1606044.4
$ws.Range("K32").Value = 1606044.4
$ws.Range("M32").Value = -1605757.4
$ws.Range("H56").Value = 69002.57000000001
$ws.Range("J56").Value = 69002.57000000001
$ws.Range("L56").Value = 69002.57000000001
$ws.Range("N56").Value = -70486.57000000001
$ws.Range("H74").Value = 36866.465
$ws.Range("I74").Value = 58299.445
$ws.Range("J74").Value = 4717
$ws.Range("K74").Value = 58299.445
$ws.Range("L74").Value = 4717
$ws.Range("M74").Value = -57425.445
$ws.Range("N74").Value = -6465
$ws.Range("H77").Value = 36866.465
$ws.Range("I77").Value = 58299.445
$ws.Range("J77").Value = 4717
$ws.Range("K77").Value = 291497.225
$ws.Range("L77").Value = 23585
$ws.Range("M77").Value = -287129.225
$ws.Range("N77").Value = -32321
$ws.Range("H132").Value = 3600.1853
$ws.Range("J132").Value = 8474.130999999999
$ws.Range("L132").Value = 25422.393
$ws.Range("N132").Value = -30482.393
$ws.Range("H133").Value = 121660
$ws.Range("J133").Value = 121660
$ws.Range("L133").Value = 121660
$ws.Range("N133").Value = -126720

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -646
$ws.Range("H94").Value = 1893.5294
$ws.Range("I94").Value = 1485.8572
$ws.Range("K94").Value = 1485.8572
$ws.Range("M94").Value = -1034.8572
$ws.Range("H99").Value = 3248796.8
$ws.Range("I99").Value = 2915.6
$ws.Range("K99").Value = 2915.6
$ws.Range("M99").Value = -1417.6
$ws.Range("H105").Value = 2887
$ws.Range("I105").Value = 1321.7222
$ws.Range("J105").Value = 4544.353
$ws.Range("K105").Value = 1321.7222
$ws.Range("L105").Value = 4544.353
$ws.Range("M105").Value = 425.2778000000001
$ws.Range("N105").Value = -8038.353
$ws.Range("H134").Value = 5465.2593
$ws.Range("I134").Value = 2076.0938
$ws.Range("K134").Value = 6228.2814
$ws.Range("M134").Value = -3693.2814

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 360.1
$ws.Range("I22").Value = 312.875
$ws.Range("K22").Value = 312.875
$ws.Range("M22").Value = 37.125
$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 2000
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2586
$ws.Range("H31").Value = 7816.102
$ws.Range("I31").Value = 4181.4644
$ws.Range("J31").Value = 12662.286
$ws.Range("K31").Value = 4181.4644
$ws.Range("L31").Value = 12662.286
$ws.Range("M31").Value = -3886.4644
$ws.Range("N31").Value = -13252.286
$ws.Range("H34").Value = 7816.102
$ws.Range("I34").Value = 4181.4644
$ws.Range("J34").Value = 12662.286
$ws.Range("K34").Value = 4181.4644
$ws.Range("L34").Value = 12662.286
$ws.Range("M34").Value = -3979.4644
$ws.Range("N34").Value = -13066.286
$ws.Range("H58").Value = 8913.311
$ws.Range("I58").Value = 3459.6667
$ws.Range("K58").Value = 3459.6667
$ws.Range("M58").Value = -3256.6667
$ws.Range("H99").Value = 4801.758
$ws.Range("I99").Value = 3437.5833
$ws.Range("K99").Value = 3437.5833
$ws.Range("M99").Value = -1939.5833
$ws.Range("H126").Value = 4801.758
$ws.Range("I126").Value = 3437.5833
$ws.Range("K126").Value = 10312.7499
$ws.Range("M126").Value = -7842.749899999999
$ws.Range("H132").Value = 5522.4585
$ws.Range("I132").Value = 3420.276
$ws.Range("K132").Value = 10260.828
$ws.Range("M132").Value = -7730.828
$ws.Range("H134").Value = 5374.593
$ws.Range("I134").Value = 2323.077
$ws.Range("J134").Value = 8208.143
$ws.Range("K134").Value = 6969.231000000001
$ws.Range("L134").Value = 24624.429
$ws.Range("M134").Value = -4434.231000000001
$ws.Range("N134").Value = -29694.429
$ws.Range("H136").Value = 8913.311
$ws.Range("I136").Value = 3459.6667
$ws.Range("K136").Value = 10379.0001
$ws.Range("M136").Value = -7829.000100000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 45458932
$ws.Range("I80").Value = 33337266
$ws.Range("K80").Value = 100011798
$ws.Range("M80").Value = -100010862
$ws.Range("H83").Value = 45458932
$ws.Range("I83").Value = 33337266
$ws.Range("K83").Value = 300035394
$ws.Range("M83").Value = -300030714
$ws.Range("H97").Value = 198.8
$ws.Range("I97").Value = 254.66667
$ws.Range("K97").Value = 764.00001
$ws.Range("M97").Value = -268.00001
$ws.Range("H131").Value = 2288.4614
$ws.Range("J131").Value = 2317.3958
$ws.Range("L131").Value = 6952.187399999999
$ws.Range("N131").Value = -17032.1874

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 60292.5
$ws.Range("J135").Value = 60292.5
$ws.Range("L135").Value = 60292.5
$ws.Range("N135").Value = -70432.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 524.2632
$ws.Range("I16").Value = 541.375
$ws.Range("K16").Value = 541.375
$ws.Range("M16").Value = -371.375
$ws.Range("H31").Value = 112918.664
$ws.Range("J31").Value = 2323.8333
$ws.Range("L31").Value = 2323.8333
$ws.Range("N31").Value = -2819.8333
$ws.Range("H122").Value = 4519.871
$ws.Range("I122").Value = 2839.9443
$ws.Range("K122").Value = 8519.832900000001
$ws.Range("M122").Value = -6069.832900000001
$ws.Range("H136").Value = 7960.643
$ws.Range("I136").Value = 2324.7222
$ws.Range("K136").Value = 6974.1666
$ws.Range("M136").Value = -4424.1666

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 78928930
$ws.Range("I14").Value = 100105750
$ws.Range("K14").Value = 100105750
$ws.Range("M14").Value = -100105582
$ws.Range("H132").Value = 19611514
$ws.Range("J132").Value = 4625.875
$ws.Range("L132").Value = 13877.625
$ws.Range("N132").Value = -18937.625

